$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46081 -> 46082, i.e. 2026-02-28 -> 2026-03-01) for every data row (2..514).
$ws.Range("C2:C514").Value = 46082
